$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Header row: replace English weekday abbreviations with Chinese weekday names ---
$ws.Range("B1").Value = "星期一"
$ws.Range("C1").Value = "星期二"
$ws.Range("D1").Value = "星期三"
$ws.Range("E1").Value = "星期四"
$ws.Range("F1").Value = "星期五"
$ws.Range("G1").Value = "星期六"
$ws.Range("H1").Value = "星期日"

# New blank header cell I1, formatted like the rest of row 1 (style s="22", same as H1)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 / D2: extend GR6017 weekly text with a new reading item (rich text) ---
$d2Text = "1.GR6017上位机衰减控制逻辑" + [char]10 + "2.GR6017 ARM Flash 读写；衰减发送" + [char]10 + "3.解决Jlink Or STLink 仿真器连接失败的问题" + [char]10 + "4.58个英语单词" + [char]10 + "5.《计算机网络自顶向下》 Cookies（初步阅读）"
$ws.Range("D2").Value = $d2Text
$redStart = $d2Text.IndexOf("初步阅读") + 1
$ws.Range("D2").Characters($redStart, 4).Font.Color = 255

# Row 2 / E2: new content describing English words + Tianao joint debugging
$e2Text = "1.79个英语单词：复习49，新增30" + [char]10 + "2.天奥联调：" + [char]10 + "    测试GR6017上位机及下位机新增的FLASH固化功能；" + [char]10 + "    以及上位机的文件保存和下载功能"
$ws.Range("E2").Value = $e2Text

# Row 2 / F2: now empty with no formatting (cell removed entirely)
$ws.Range("F2").Clear()

# --- Row 1 height / D column width updates ---
$ws.Rows.Item(2).RowHeight = 132.75
$ws.Columns.Item(4).ColumnWidth = 21.8

# --- Selection moves to E4 ---
$ws.Range("E4").Select()
